$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Snapshot the existing cell-format combinations (font/fill/border/numFmt)
#    we will need for the new 6-column layout, into a scratch row far below
#    the table, BEFORE any values/formats are overwritten.
# ---------------------------------------------------------------------------
# xf7  (A1 header-left)   <- A1
# xf1  (B1/C1 header)     <- B1
# xf2  (D1/E1 header)     <- D1
# xf8  (E1 header-right)  <- E1
# xf9  (A2:A4 body-left)  <- B2
# xf10 (F2:F4 body-right) <- A2
# xf11 (A5 footer-left)   <- E2
# xf12 (B5:E5 footer)     <- A5
# xf13 (F5 footer-right)  <- B5
# xf14 (B2:E4 body)       <- E5

$ws.Range("B1").Copy()
$ws.Range("A200").PasteSpecial(-4122)

$ws.Range("D1").Copy()
$ws.Range("B200").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("C200").PasteSpecial(-4122)

$ws.Range("E1").Copy()
$ws.Range("D200").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("E200").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("F200").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("G200").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("H200").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("I200").PasteSpecial(-4122)

$ws.Range("E5").Copy()
$ws.Range("J200").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Apply the snapshotted formats onto the final 6-column cell groups.
# ---------------------------------------------------------------------------
$ws.Range("C200").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A200").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)

$ws.Range("B200").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

$ws.Range("D200").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("E200").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

$ws.Range("J200").Copy()
$ws.Range("B2:E4").PasteSpecial(-4122)

$ws.Range("F200").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122)

$ws.Range("G200").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("H200").Copy()
$ws.Range("B5:E5").PasteSpecial(-4122)

$ws.Range("I200").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Remove the scratch cells (contents + formats), so only A1:F5 remain used.
# ---------------------------------------------------------------------------
$ws.Range("A200:J200").Clear()

# ---------------------------------------------------------------------------
# 4) Write the new cell text / values (order matters for some NumberFormats
#    set above; text values go through shared strings automatically).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "位号/Designator"
$ws.Range("B1").Value = "型号/Comment"
$ws.Range("C1").Value = "封装/Footprint"
$ws.Range("D1").Value = "数量/Quantity"
$ws.Range("E1").Value = "备注/Description"
$ws.Range("F1").Value = "嘉立创元件编号"

$ws.Range("A2").Value = "J1"
$ws.Range("B2").Value = "USB_C_Receptacle_USB2.0"
$ws.Range("C2").Value = "USB_C_Receptacle_HRO_TYPE-C-31-M-12"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "16Pin Type-C母座"
$ws.Range("F2").Value = "C709357"

$ws.Range("A3").Value = "J2"
$ws.Range("B3").Value = "JST_SH_SM04B-SRSS"
$ws.Range("C3").Value = "JST_SH_SM04B-SRSS-TB_1x04-1MP_P1.00mm_Horizontal"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "SH 4Pin端子插座"
$ws.Range("F3").Value = "C2764183"

$ws.Range("A4").Value = "R1, R2"
$ws.Range("B4").Value = "5.1kΩ"
$ws.Range("C4").Value = "R_0603"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "5%贴片电阻"

$ws.Range("A5").Value = "U1"
$ws.Range("B5").Value = "USBLC6-2P6"
$ws.Range("C5").Value = "SOT-666"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "ESD芯片"
$ws.Range("F5").Value = "C2827693"

# ---------------------------------------------------------------------------
# 5) Column widths for the new layout (A..F).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.81640625
$ws.Columns.Item(2).ColumnWidth = 35.453125
$ws.Columns.Item(3).ColumnWidth = 48.6328125
$ws.Columns.Item(4).ColumnWidth = 18.54296875
$ws.Columns.Item(5).ColumnWidth = 31.26953125
$ws.Columns.Item(6).ColumnWidth = 39.81640625

# ---------------------------------------------------------------------------
# 6) Row heights.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 24
$ws.Rows.Item(5).RowHeight = 17

# ---------------------------------------------------------------------------
# 7) Selection matches the post-edit state.
# ---------------------------------------------------------------------------
$ws.Range("D9").Select()
